# Actualización automática 2025-07-08 16:45:08
#
# Inserts a new data row (OFICINA-CATAECSA / VEHINVER SA, all-zero metrics)
# at row 250 in both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets,
# pushing the existing rows 250-282 down to 251-283. The trailing
# summary row (now row 283) is updated so its "N de 280" style labels
# read "N de 281" to reflect the new total row count.

$wb = $excel.ActiveWorkbook

# ---- Sheet "VENTAS POR GRUPO" (columns A:R) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(250).Insert()

$ws1.Range("A250").Value = "OFICINA-CATAECSA"
$ws1.Range("B250").Value = "VEHINVER SA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(250, $c).Value = 0
}

# Fix the "N de 280" -> "N de 281" summary labels on the (now shifted)
# totals row, which moved from row 282 to row 283.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(283, $c)
    $cur = $cell.Value()
    $cell.Value = $cur.Replace("de 280", "de 281")
}

# ---- Sheet "VENTA MENSUAL" (columns A:G) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(250).Insert()

$ws2.Range("A250").Value = "OFICINA-CATAECSA"
$ws2.Range("B250").Value = "VEHINVER SA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(250, $c).Value = 0
}
